# Lab01_ReviewReport.xlsx - "Added architectural design phase report"
#
# This script reproduces the content-level edits described by the commit:
#   - Requirements Phase Defects sheet: mark the three existing requirement
#     defects (R01/R02/R03) with a "1" count in column D, and move the
#     active-cell selection to D13.
#   - Architect. Design Phase Defects sheet: fill in the first three defect
#     rows (A01/A02/A03) with their counts and comments, taller rows to fit
#     the wrapped text, and move the active-cell selection to E15.

$wb = $excel.ActiveWorkbook

$reqSheet  = $wb.Worksheets.Item("Requirements Phase Defects")
$archSheet = $wb.Worksheets.Item("Architect. Design Phase Defects")

# --- Requirements Phase Defects: tick the "count" column for the 3 existing rows ---
$reqSheet.Cells.Item(10, 4).Value = 1
$reqSheet.Cells.Item(11, 4).Value = 1
$reqSheet.Cells.Item(12, 4).Value = 1

# --- Architect. Design Phase Defects: add the new architectural defects ---
$archSheet.Rows.Item(10).RowHeight = 30
$archSheet.Rows.Item(11).RowHeight = 30
$archSheet.Rows.Item(12).RowHeight = 30

# "Crt. item" column (matches the shared-string build order of the original edit)
$archSheet.Cells.Item(10, 3).Value = "A01"
$archSheet.Cells.Item(11, 3).Value = "A02"
$archSheet.Cells.Item(12, 3).Value = "A03"

# "Comments/ improvements" column
$archSheet.Cells.Item(10, 5).Value = "Architecture diagram is chaotic and not comprehensive"
$archSheet.Cells.Item(12, 5).Value = "The architecure does not define the entities used for the application"
$archSheet.Cells.Item(11, 5).Value = "Useless packages/classes and layering issues are present"

# "Doc. page/line" (count) column
$archSheet.Cells.Item(10, 4).Value = 1
$archSheet.Cells.Item(11, 4).Value = 1
$archSheet.Cells.Item(12, 4).Value = 1

# --- window / selection bookkeeping ---
# Move the window position the same way the original author's Excel session
# did when it re-saved the file.
$win = $excel.ActiveWindow
$win.Top = 1500

# Update the remembered cell selection on each affected sheet.
$archSheet.Range("E15").Select() | Out-Null
$reqSheet.Activate() | Out-Null
$reqSheet.Range("D13").Select() | Out-Null
